# Weekly update: a new week of "Cilantro" price data (Mercado Mayorista Lo
# Valledor de Santiago) was published, so two new records (rows) are
# inserted right after the existing row 553, pushing the rest of the table
# down by two rows. The two brand-new rows carry the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 554-555; everything from the old row 554 onward
# shifts down to 556 onward (old 586-587 becomes the new 588-589).
$ws.Rows("554:555").Insert()

# --- New row 554 -----------------------------------------------------
$ws.Range("A554").Value = 6
$ws.Range("B554").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C554").Value = "Metropolitana"
$ws.Range("D554").Value = 44516
$ws.Range("E554").Value = 13
$ws.Range("F554").Value = 100112040
$ws.Range("G554").Value = "Cilantro"
$ws.Range("H554").Value = "Sin especificar"
$ws.Range("I554").Value = "Primera"
$ws.Range("J554").Value = 440
$ws.Range("K554").Value = 9000
$ws.Range("L554").Value = 10000
$ws.Range("M554").Value = 9341
$ws.Range("N554").Value = "$/caja 36 atados"
$ws.Range("O554").Value = "Región Metropolitana"
$ws.Range("P554").Value = 259
$ws.Range("Q554").Value = 36
$ws.Range("R554").Value = "Hortaliza"

# --- New row 555 -----------------------------------------------------
$ws.Range("A555").Value = 6
$ws.Range("B555").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C555").Value = "Metropolitana"
$ws.Range("D555").Value = 44516
$ws.Range("E555").Value = 13
$ws.Range("F555").Value = 100112040
$ws.Range("G555").Value = "Cilantro"
$ws.Range("H555").Value = "Sin especificar"
$ws.Range("I555").Value = "Primera"
$ws.Range("J555").Value = 370
$ws.Range("K555").Value = 14000
$ws.Range("L555").Value = 15000
$ws.Range("M555").Value = 14405
$ws.Range("N555").Value = "$/docena de atados"
$ws.Range("O555").Value = "Región Metropolitana"
$ws.Range("P555").Value = 4802
$ws.Range("Q555").Value = 3
$ws.Range("R555").Value = "Hortaliza"
